$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.148.52'
$ws.Range("E2").Value = '  +1.18%  '
$ws.Range("D3").Value = '2.055.69'
$ws.Range("E3").Value = '  -3.58%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '''248.76'
$ws.Range("E5").Value = '  -2.72%  '
$ws.Range("E6").Value = '  -2.20%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '''54.83'
$ws.Range("E8").Value = '  +16.34%  '
$ws.Range("D9").Value = '''61.87'
$ws.Range("E9").Value = '  +3.21%  '
$ws.Range("D10").Value = '''0.376'
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("E11").Value = '  +5.70%  '
$ws.Range("E12").Value = '  +5.11%  '
$ws.Range("E13").Value = '  +4.89%  '
$ws.Range("D14").Value = '2.352.88'
$ws.Range("E14").Value = '  -3.62%  '
$ws.Range("D15").Value = '''0.815'
$ws.Range("E15").Value = '  -3.35%  '
$ws.Range("D16").Value = '''5.22'
$ws.Range("E16").Value = '  +1.76%  '
$ws.Range("D17").Value = '2.053.32'
$ws.Range("E17").Value = '  -3.59%  '
$ws.Range("D18").Value = '37.114.31'
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").Value = '''72.26'
$ws.Range("E19").Value = '  -1.88%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0901'
$ws.Range("E20").Value = '  +7.58%  '
$ws.Range("D21").Value = '''14.26'
$ws.Range("E21").Value = '  +6.80%  '
$ws.Range("D22").Value = '''5.32'
$ws.Range("E22").Value = '  +2.10%  '
$ws.Range("D23").Value = '''236.74'
$ws.Range("E23").Value = '  -2.02%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("E25").Value = '  -2.74%  '
$ws.Range("D26").Value = '''170.06'
$ws.Range("E26").Value = '  -1.12%  '
$ws.Range("E27").Value = '  -2.27%  '
$ws.Range("D28").Value = '''20.19'
$ws.Range("E28").Value = '  -7.54%  '
$ws.Range("E29").Value = '  -2.64%  '
$ws.Range("E30").Value = '  -0.59%  '
$ws.Range("E31").Value = '  +0.90%  '
$ws.Range("E32").Value = '  +11.40%  '
$ws.Range("E33").Value = '  +3.82%  '
$ws.Range("E34").Value = '  +3.41%  '
$ws.Range("D35").Value = '''0.0882'
$ws.Range("E35").Value = '  -8.10%  '
$ws.Range("D36").Value = '''1.00'
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  -5.04%  '
$ws.Range("D38").Value = '''1.74'
$ws.Range("E38").Value = '  -8.21%  '
$ws.Range("E39").Value = '  -0.58%  '
$ws.Range("E40").Value = '  +22.73%  '
$ws.Range("D41").Value = '''18.25'
$ws.Range("E41").Value = '  +12.96%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.0223'
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("B43").Value = 'Gas'
$ws.Range("C43").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D43").Value = '''15.35'
$ws.Range("E43").Value = '  -46.03%  '
$ws.Range("E44").Value = '  -5.49%  '
$ws.Range("D45").Value = '''95.96'
$ws.Range("E45").Value = '  -3.15%  '
$ws.Range("D46").Value = '''2.78'
$ws.Range("E46").Value = '  -1.08%  '
$ws.Range("D47").Value = '''4.18'
$ws.Range("E47").Value = '  +36.41%  '
$ws.Range("D48").Value = '1.295.60'
$ws.Range("E48").Value = '  -4.75%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''2.37'
$ws.Range("E49").Value = '  +3.16%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = '''2.92'
$ws.Range("E50").Value = '  +2.84%  '
$ws.Range("D51").Value = '''6.77'
$ws.Range("E51").Value = '  -7.14%  '
